$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.722.76"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").Value = "1.880.05"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.53%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4707"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3937"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.54"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08060"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.026"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.67%  "

$ws.Range("D13").Value = "1.887.14"
$ws.Range("E13").Value = "  +1.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.975"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.142"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001050"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.10%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06707"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.543"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.72%  "

$ws.Range("D23").Value = "27.735.32"
$ws.Range("E23").Value = "  +1.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.316"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.35%  "

$ws.Range("D26").Value = "2.110.63"
$ws.Range("E26").Value = "  +1.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.107"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.594"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9852"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09482"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.454"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.615"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.359"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06144"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02266"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.231"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.129"
$ws.Range("D40").Style = "Normal"

$ws.Range("E41").Value = "  +1.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1901"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.257"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5704"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.52%  "

$ws.Range("E46").Value = "  +2.05%  "

$ws.Range("E47").Value = "  +1.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.396"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06908"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000301"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.17%  "
